# Apply the CodeSystem-wh-payer-subscriber-relationship metadata update:
#  - Version bump 5.0.0 -> 6.0.0
#  - Date bump
#  - Publisher gets a real value ("Alvearie Team")
#  - The old duplicate "Contact" / "No display for ContactDetail" row is
#    replaced with a new "Jurisdiction" / "United States of America" row
#    and the second (duplicate) Contact row is removed entirely
#  - "Case Sensitive" gets a value of "true"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: refreshed timestamp
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value
$ws.Range("B9").Value = "Alvearie Team"

# First "Contact" row becomes "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# The second (duplicate) "Contact" row (row 11) is removed entirely,
# shifting Description/Purpose/Copyright/Case Sensitive/... up by one row
$ws.Rows.Item(11).Delete()

# "Case Sensitive" (now row 14 after the deletion above) gets value "true".
# A plain $cell.Value = "true" gets auto-coerced to a native boolean by
# Excel's "smart" literal parsing (just like typing it into the grid), so
# instead put a formula that evaluates to the text "true" and then convert
# that formula to its literal value via copy / paste-special-values - this
# keeps the cell a genuine text cell (t="s") instead of t="b".
$c = $ws.Range("B14")
$c.Formula = "=""true"""
$c.Copy()
$c.PasteSpecial(-4163)
